$wb = $excel.ActiveWorkbook

# --- Sheet "Schedule": update Cost ($) and Unit Cost ($/ML) for rows 3 and 4 ---
$schedule = $wb.Worksheets.Item("Schedule")
$schedule.Range("E3").Value = 765.6930742499999
$schedule.Range("F3").Value = 28.93775790816326
$schedule.Range("E4").Value = 383.70008625
$schedule.Range("F4").Value = 11.27866214726632

# --- Sheet "Detailed": update Price column (B) and Type column (C) ---
$detailed = $wb.Worksheets.Item("Detailed")

$priceUpdates = @{
    41 = 51.4753
    42 = 59.36649
    43 = 57.09
    44 = 57.03884
    45 = 57.09
    46 = 57.06
    47 = 57.09
    52 = 56.98
    53 = 56.98
    54 = 41.41096
    61 = 57.12044
    62 = 57.36498
    64 = 26.87623
    65 = 34.31377
    66 = 25.84442
    70 = 5.34096
    71 = 23.48462
    72 = 22.82458
    73 = 26.22731
    74 = 26.99167
    77 = 36.06
    79 = 29.51494
    80 = 27.88341
    81 = -2.10336
    82 = -2.94351
    85 = -4.04163
    86 = 12.0408
    87 = 33.06721
    89 = 57.06
    90 = 57.06
    91 = 56.98
    92 = 47.19381
    93 = 54.85437
    94 = 50.62881
    97 = 52.12057
}

foreach ($row in $priceUpdates.Keys) {
    $detailed.Cells.Item($row, 2).Value = $priceUpdates[$row]
}

# Type (column C) changes: rows 43 and 44 go from "forecast" to "historical"
$detailed.Cells.Item(43, 3).Value = "historical"
$detailed.Cells.Item(44, 3).Value = "historical"
